# Exam workbook update:
#  - rename the "result" column header to "ref_answer"
#  - rephrase a handful of exam question prompts (column B) into
#    alternate wording while keeping the same answer choices/answers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "ref_answer"

$ws.Range("B6").Value = "ทุนพัฒนาบุคลากร อยูในประกาศ หลักเกณฑ์การจ่ายทุนพัฒนาและส่งเสริมศักยภาพผู้เรียน หรือไม่"
$ws.Range("B17").Value = "ค่าใช้จ่ายในการเดินทางรายวันในกรุงเทพมหานคร ครอบคลุมทุกประเภททุนหรือไม่"
$ws.Range("B26").Value = "ความรับผิดชอบของผู้รับทุนที่ถูกระบุไว้ในประกาศ(หลักเกณฑ์การจ่ายทุนพัฒนาและส่งเสริมศักยภาพผู้เรียน) คืออะไร"
$ws.Range("B31").Value = "ประเภทของทุนช่วยเหลือการศึกษาในประกาศ หลักเกณฑ์การจ่ายทุนช่วยเหลือการศึกษาสำหรับผู้ขาดแคลนทุนทรัพย์ คือทุนประเภทอะไร"
$ws.Range("B40").Value = "ทุนค่าใช้จ่ายรายเดือนในกรณีที่ผู้เรียนได้รับทุนอื่นร่วมด้วยต้องทำยังไง"
$ws.Range("B42").Value = "อะไรไม่ใช่หน้าที่ของอธิการบดีในประกาศนี้(หลักเกณฑ์การจ่ายทุนช่วยเหลือการศึกษาสำหรับผู้ขาดแคลนทุนทรัพย์)?"
$ws.Range("B43").Value = "อะไรคือหลักเกณฑ์สำคัญในการขอรับทุนช่วยเหลือการศึกษาสำหรับผู้ขาดแคลนทุนทรัพย์?"
$ws.Range("B64").Value = "ทุนส่งเสริมศักยภาพเหมาะสำหรับผู้สมัครประเภทใดเป็นหลัก?"
$ws.Range("B89").Value = "อะไรคือเงื่อนไขเพิ่มเติมสำหรับผู้ที่ได้รับทุนประเภท ง?"

[void]$ws.Range("F3").Select()
